$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Regresión Logística gray"
$ws.Range("B2").Value = 0.2377205781933585
$ws.Range("C2").Value = 0.2369306836345303
$ws.Range("D2").Value = 0.2344130759382817
$ws.Range("E2").Value = 0.2454802259887006
$ws.Range("F2").Value = 0.5362792635039743
